$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(53, 1).Value = 52
$ws.Cells.Item(53, 2).Value = 43882
$ws.Cells.Item(53, 3).Value = "Surah Al Nisa, 108 - 119"
$ws.Cells.Item(53, 4).Value = "Test content"
$ws.Cells.Item(53, 5).Value = "Qasim Ali"
$ws.Cells.Item(53, 6).Value = "Reacting to criticism, Importance of purpose in life, Private conversations, Living a good life"

$ws.Range("B53").NumberFormat = "d-mmm-yy"

$ws.Rows.Item(53).RowHeight = 409.6
